# changed demo module name
# - Rename the "Demo" section title to "Code"
# - Collapse a couple of split runs back into single runs (no visible text
#   change, just fewer runs), matching the author's re-save of those slides.

$p = $ppt.ActivePresentation

# Slide 10: Title "Demo" -> "Code"
$slide10 = $p.Slides.Item(10)
$slide10.Shapes.Item(1).TextFrame.TextRange.Text = "Code"

# Slide 2: merge the two trailing runs of the subtitle into a single run.
# "owershell is a journey. Each and every one of " + "us is at a different
# point in that journey. Knowing the path forward can often be the most
# challenging part if you don't have a guide. "
$slide2 = $p.Slides.Item(2)
$subtitle2 = $slide2.Shapes.Item(2).TextFrame.TextRange
$merge2 = $subtitle2.Characters(11, 181)
$merge2.Text = "owershell is a journey. Each and every one of us is at a different point in that journey. Knowing the path forward can often be the most challenging part if you don’t have a guide. "

# Slide 3: merge the leading " " run and "The " run into one " The " run.
$slide3 = $p.Slides.Item(3)
$title3 = $slide3.Shapes.Item(1).TextFrame.TextRange
$merge3 = $title3.Characters(1, 5)
$merge3.Text = " The "

# Slide 4: same merge as slide 3.
$slide4 = $p.Slides.Item(4)
$title4 = $slide4.Shapes.Item(1).TextFrame.TextRange
$merge4 = $title4.Characters(1, 5)
$merge4.Text = " The "
